$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: switch classification from iaest-dimension: to iaest-measure: for columns A, D, J
$ws.Range("A2").Value = "iaest-measure:temporalidad"
$ws.Range("D2").Value = "iaest-measure:mes-nombre"
$ws.Range("J2").Value = "iaest-measure:sexo"

# Row 3: switch "dim" to "medida" for columns A, D, J
$ws.Range("A3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: switch "skos:Concept" to "xsd:int" for columns A, D, J
$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: remove mapping file references for columns A, D, J (no longer dimensions)
$ws.Range("A5").Clear()
$ws.Range("D5").Clear()
$ws.Range("J5").Clear()
